# Release Burndown Chart.xlsx — update "Tasks to Complete" burndown figures
# for sprint points 0 and 1 (row 11, cols D and E), and move the active
# selection to E14 (matches the recorded sheet-view state after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = 106
$ws.Range("E11").Value = 56

$ws.Range("E14").Select()
